# Insert a new data row at row 121 (shifts existing rows 121-203 down to 122-204)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("121").Insert()

$ws.Range("A121").Value = 9
$ws.Range("B121").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C121").Value = "Metropolitana"
$ws.Range("D121").Value = 44673
$ws.Range("E121").Value = 13
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100101
$ws.Range("H121").Value = "Berries"
$ws.Range("I121").Value = 100101001
$ws.Range("J121").Value = "Arándano (blue)"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 329
$ws.Range("N121").Value = 5000
$ws.Range("O121").Value = 5000
$ws.Range("P121").Value = 5000
$ws.Range("Q121").Value = "$/bandeja 2 kilos"
$ws.Range("R121").Value = "Provincia de Linares"
$ws.Range("S121").Value = 2500
$ws.Range("T121").Value = 2
